# WW101-Binder-Cover.pptx edit
#  - refresh the "datetimeFigureOut" date placeholder on the slide master
#    and every slide layout: 1/31/2017 -> 2/1/2017
#  - slide 1, "TextBox 5" (author list): tidy up the existing runs and
#    add Mike Noel + Vikram Ramanna to the list

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders (slide master + all custom layouts)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "1/31/2017") {
                $shp.TextFrame.TextRange.Text = "2/1/2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

# ---------------------------------------------------------------------
# 2) Slide 1 author-credits textbox
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$credits = $null
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "*Greg Landry*") {
        $credits = $shp
    }
}

$tr = $credits.TextFrame.TextRange

# "Greg Landry" -> split into "Greg " / "Landry" runs, then add "Mike Noel"
$greg = $tr.Paragraphs(2, 1)
$gregSplit = $greg.Characters(1, 5)
$gregSplit.Text = "Greg_"
$gregSplit = $greg.Characters(1, 5)
$gregSplit.Text = "Greg "
$greg = $tr.Paragraphs(2, 1)
$greg.InsertAfter("`rMike Noel") | Out-Null

# "James Dougherty" -> split into "James " / "Dougherty" runs, then add "Vikram Ramanna"
$james = $tr.Paragraphs(4, 1)
$jamesSplit = $james.Characters(1, 6)
$jamesSplit.Text = "James_"
$jamesSplit = $james.Characters(1, 6)
$jamesSplit.Text = "James "
$james = $tr.Paragraphs(4, 1)
$james.InsertAfter("`rVikram Ramanna") | Out-Null

# The textbox auto-fits its height to the (now longer) text; pin the
# exact EMU value the real deck ends up with (1938992 EMU == 6858000 x
# extent's companion cy) to avoid a sub-EMU float round-trip diff.
$credits.Height = (1938992 + 0.5) / 12700
